# Työaikakirjanpito update: add a new logged day (row 42) describing
# Adobe XD hifi-version work + logo redesign, and move the active
# selection to F32 (matches the author's next working cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42: new time-log entry ------------------------------------------
# Päivä (date serial 44070 = 27.8.2020)
$ws.Range("A42").Value = 44070
# Aika (h)
$ws.Range("B42").Value = 6.5
# Mitä tein - new shared string, wrapped like the other description cells
$ws.Range("C42").WrapText = $true
$ws.Range("C42").Value = "Adobe XD:llä hifi version tekemistä ja logon uusiminen"

# Row grows taller to fit the wrapped description, same as the other
# multi-line rows in the sheet.
$ws.Rows.Item(42).RowHeight = 30

# --- Move the active selection -------------------------------------------
$ws.Activate()
$ws.Range("F32").Select()
